$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the split "To use the pattern..." run (separated by the hidden
#    _GoBack bookmark) into a single contiguous sentence ending
#    "...the actual data in the client."
# ---------------------------------------------------------------------------

$mergeIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("To use the pattern: you can use static data members")) {
        $mergeIndex = $i
        break
    }
}

if ($mergeIndex -eq -1) {
    throw "Could not locate the Flyweight 'To use the pattern' paragraph"
}

$mergeParagraph = $d.Paragraphs($mergeIndex)
$mergeRange = $mergeParagraph.Range
$tail = " the actual data in the client."
$tailRange = $d.Range($mergeRange.End - ($tail.Length + 1), $mergeRange.End)
$tailRange.Text = ""

$mergeParagraph2 = $d.Paragraphs($mergeIndex)
$mergeParagraph2.Range.InsertAfter($tail)

# ---------------------------------------------------------------------------
# 2. Find the "Proxy" heading and the placeholder bullet paragraph (">")
#    directly beneath it, then turn it into the four real bullet points
#    describing the Proxy pattern.
# ---------------------------------------------------------------------------

$proxyHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Style.NameLocal -eq "C_Heading" -and $pp.Range.Text.Trim() -eq "Proxy") {
        $proxyHeadingIndex = $i
        break
    }
}

if ($proxyHeadingIndex -eq -1) {
    throw "Could not locate the Proxy heading paragraph"
}

$bulletIndex = $proxyHeadingIndex + 1
$bulletParagraph = $d.Paragraphs($bulletIndex)

if ($bulletParagraph.Range.Text.Trim() -ne ">") {
    throw "Proxy placeholder bullet paragraph not found where expected"
}

# First bullet: replace the ">" placeholder text in place.
$bulletParagraph.Range.Text = "The proxy design pattern functions as an interface to a particular resource. That resource may be remote, expensive to constructor, or may require logging or some other added functionality."

# Second bullet.
$d.Paragraphs($bulletIndex).Range.InsertParagraphAfter()
$bullet2 = $d.Paragraphs($bulletIndex + 1)
$bullet2.Range.Text = "Example: Smart pointers are proxies since for the most part they work like raw pointers. They are enhanced however since they offer additional functionality that raw pointers do not allow."

# Third bullet.
$d.Paragraphs($bulletIndex + 1).Range.InsertParagraphAfter()
$bullet3 = $d.Paragraphs($bulletIndex + 2)
$bullet3.Range.Text = "Example: Children on the same level in the class hierarchy are proxies for one another. By replacing one type for another in a variable declaration, the code compiles because the interface is unchanged. "

# Fourth bullet.
$d.Paragraphs($bulletIndex + 2).Range.InsertParagraphAfter()
$bullet4 = $d.Paragraphs($bulletIndex + 3)
$bullet4.Range.Text = "To use the pattern: create another class that has the exact same interface as the target class that it is proxying. You can use inheritance, or composition. The interface of one class being the same as another is the main focus, how you construct or declare the object is irrelevant."

Write-Output "Proxy section populated; document now has $($d.Paragraphs.Count) paragraphs."
